# Refresh cryptocurrency price (D) and 1h volume-change (E) columns
# with the latest values from the scheduled scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.823.91"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "1.734.33"
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'230.37"
$ws.Range("E5").Value = "  -2.85%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.5224"
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").Value = "'0.2754"
$ws.Range("E8").Value = "  +1.52%  "
$ws.Range("D9").Value = "'39.27"
$ws.Range("E9").Value = "  -3.05%  "
$ws.Range("D10").Value = "'0.06133"
$ws.Range("E10").Value = "  -1.09%  "
$ws.Range("D11").Value = "1.741.55"
$ws.Range("E11").Value = "  -1.73%  "
$ws.Range("D12").Value = "'0.07040"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").Value = "'15.01"
$ws.Range("E13").Value = "  -4.25%  "
$ws.Range("D14").Value = "'0.6364"
$ws.Range("E14").Value = "  -2.52%  "
$ws.Range("D15").Value = "'4.524"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").Value = "'76.63"
$ws.Range("E16").Value = "  -1.97%  "
$ws.Range("D17").Value = "'0.9998"
$ws.Range("D18").Value = "'0.9997"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "25.823.21"
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").Value = "'11.47"
$ws.Range("E20").Value = "  -1.76%  "
$ws.Range("D21").Value = "'0.000006635"
$ws.Range("E21").Value = "  -0.99%  "
$ws.Range("D22").Value = "1.954.91"
$ws.Range("E22").Value = "  -1.74%  "
$ws.Range("D23").Value = "'4.188"
$ws.Range("E23").Value = "  +2.38%  "
$ws.Range("D24").Value = "'8.784"
$ws.Range("E24").Value = "  +4.87%  "
$ws.Range("D25").Value = "'5.158"
$ws.Range("E25").Value = "  -0.38%  "
$ws.Range("D26").Value = "'139.58"
$ws.Range("E26").Value = "  +1.60%  "
$ws.Range("D27").Value = "'1.503"
$ws.Range("E27").Value = "  +1.66%  "
$ws.Range("D28").Value = "'15.00"
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("D29").Value = "'1.781"
$ws.Range("E29").Value = "  -2.27%  "
$ws.Range("D30").Value = "'102.05"
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("D31").Value = "'0.08295"
$ws.Range("E31").Value = "  -1.46%  "
$ws.Range("D32").Value = "'3.713"
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").Value = "'3.498"
$ws.Range("E33").Value = "  +2.22%  "
$ws.Range("E34").Value = "  +1.16%  "
$ws.Range("D35").Value = "'2.604"
$ws.Range("E35").Value = "  -1.80%  "
$ws.Range("D36").Value = "'0.9730"
$ws.Range("E36").Value = "  -2.49%  "
$ws.Range("D37").Value = "'0.6165"
$ws.Range("E37").Value = "  +1.40%  "
$ws.Range("D38").Value = "'2.666"
$ws.Range("E38").Value = "  -2.57%  "
$ws.Range("E39").Value = "  -0.33%  "
$ws.Range("D40").Value = "'0.9996"
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("D41").Value = "'1.900"
$ws.Range("E41").Value = "  -2.29%  "
$ws.Range("E42").Value = "  -2.87%  "
$ws.Range("D43").Value = "'0.3824"
$ws.Range("E43").Value = "  -1.26%  "
$ws.Range("D44").Value = "'5.015"
$ws.Range("E44").Value = "  +1.69%  "
$ws.Range("D45").Value = "'0.7205"
$ws.Range("E45").Value = "  -4.29%  "
$ws.Range("D46").Value = "'0.05342"
$ws.Range("E46").Value = "  -2.80%  "
$ws.Range("D47").Value = "'0.1126"
$ws.Range("E47").Value = "  +0.75%  "
$ws.Range("D48").Value = "'6.176"
$ws.Range("E48").Value = "  +1.31%  "
$ws.Range("D49").Value = "'53.15"
$ws.Range("E49").Value = "  +0.82%  "
$ws.Range("D50").Value = "'29.95"
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("D51").Value = "'7.625"
$ws.Range("E51").Value = "  +2.77%  "
